$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.075.25"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").Value = "3.542.15"
$ws.Range("E3").Value = "  -0.40%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.10"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.60"
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("D7").Value = "3.543.01"
$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.77"
$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.406"
$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("D13").Value = "4.143.14"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000200"
$ws.Range("E14").Value = "  -3.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.89"
$ws.Range("E15").Value = "  -4.12%  "

$ws.Range("D16").Value = "3.555.69"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  +1.70%  "

$ws.Range("D18").Value = "66.063.64"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.00"
$ws.Range("E19").Value = "  -4.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.63"
$ws.Range("E21").Value = "  -1.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "416.37"
$ws.Range("E22").Value = "  -3.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.599"
$ws.Range("E23").Value = "  -1.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.75"
$ws.Range("E24").Value = "  -2.21%  "

$ws.Range("D25").Value = "3.684.50"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("E27").Value = "  -3.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.04"
$ws.Range("E28").Value = "  -1.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.46"
$ws.Range("E29").Value = "  -1.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.81"
$ws.Range("E30").Value = "  -1.75%  "

$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("D32").Value = "3.537.35"
$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("E33").Value = "  +2.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.40"
$ws.Range("E34").Value = "  -4.06%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.55"
$ws.Range("E36").Value = "  -3.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.30"
$ws.Range("E37").Value = "  -10.18%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "173.74"
$ws.Range("E38").Value = "  -1.30%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.60"
$ws.Range("E39").Value = "  -7.44%  "

$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.25"
$ws.Range("E40").Value = "  -6.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0818"
$ws.Range("E41").Value = "  -3.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.06"
$ws.Range("E42").Value = "  -2.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.859"
$ws.Range("E43").Value = "  -3.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.56"
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.80"
$ws.Range("E45").Value = "  -6.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").Value = "  -4.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.07"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.67"
$ws.Range("E49").Value = "  -2.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.10"
$ws.Range("E50").Value = "  -8.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.16"
$ws.Range("E51").Value = "  -8.28%  "
